$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.5
$ws.Range("C2").Value = 5.344827651977539
$ws.Range("D2").Value = 11.11936092376709
$ws.Range("E2").Value = 34.28571701049805
